# Auto-generated Excel COM-interop script
# Applies the cryptos.xlsx price/volume update described in the commit:
# "Updated cryptos list on Fri Nov 29 21:26:50 UTC 2024 with GitHub Actions"
#
# The sheet stores every value as text (coin name, coinranking.com link,
# a '.'-punctuated price string and a padded/percent-sign volume string).
# Column D prices are things like '44.09', '0.170' or '97.383.86' -- some
# are plain decimals that Excel would otherwise silently re-interpret as
# IEEE-754 numbers (losing trailing zeros / gaining float noise such as
# '44.090000000000003'), so each Price cell being rewritten is switched to
# Text format ('@') right before its new value is written, keeping the
# on-disk text byte-for-byte what the source feed produced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '97.126.49'
$ws.Range("E2").Value = '  +2.11%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.579.39'
$ws.Range("E3").Value = '  +0.27%  '

# Row 4
$ws.Range("E4").Value = '  +0.00%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.77'
$ws.Range("E5").Value = '  +2.77%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.75'
$ws.Range("E6").Value = '  +18.02%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '651.18'
$ws.Range("E7").Value = '  -0.52%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.428'
$ws.Range("E8").Value = '  +6.98%  '

# Row 9
$ws.Range("E9").Value = '  -0.08%  '

# Row 10
$ws.Range("E10").Value = '  +3.25%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.575.54'
$ws.Range("E11").Value = '  +0.19%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '44.09'
$ws.Range("E12").Value = '  +3.99%  '

# Row 13
$ws.Range("E13").Value = '  +0.78%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.46'
$ws.Range("E14").Value = '  +0.12%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.247.50'
$ws.Range("E15").Value = '  +0.33%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '96.879.65'
$ws.Range("E16").Value = '  +1.93%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000260'
$ws.Range("E17").Value = '  +2.76%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.63'
$ws.Range("E18").Value = '  +1.05%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.589.17'
$ws.Range("E19").Value = '  +0.52%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.55'
$ws.Range("E20").Value = '  -1.51%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '18.03'
$ws.Range("E21").Value = '  +1.08%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.531'
$ws.Range("E22").Value = '  +10.50%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '516.73'
$ws.Range("E23").Value = '  +1.59%  '

# Row 24
$ws.Range("E24").Value = '  +0.72%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000207'
$ws.Range("E25").Value = '  +5.47%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.89'
$ws.Range("E26").Value = '  -0.36%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '103.38'
$ws.Range("E27").Value = '  +8.53%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '13.09'
$ws.Range("E28").Value = '  +2.97%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.773.45'
$ws.Range("E29").Value = '  +0.28%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.170'
$ws.Range("E30").Value = '  +17.51%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.97'
$ws.Range("E31").Value = '  +3.60%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.97'
$ws.Range("E32").Value = '  -1.94%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  -0.08%  '

# Row 34
$ws.Range("E34").Value = '  +6.11%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  +0.04%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '31.75'
$ws.Range("E36").Value = '  -0.46%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.574'
$ws.Range("E37").Value = '  +2.39%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '617.60'
$ws.Range("E38").Value = '  +2.09%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.67'
$ws.Range("E39").Value = '  +0.64%  '

# Row 40
$ws.Range("E40").Value = '  -5.30%  '

# Row 41
$ws.Range("E41").Value = '  +2.13%  '

# Row 42
$ws.Range("E42").Value = '  +5.48%  '

# Row 43
$ws.Range("E43").Value = '  -0.06%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.925'
$ws.Range("E44").Value = '  +2.11%  '

# Row 45
$ws.Range("B45").Value = 'Algorand'
$ws.Range("C45").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.438'
$ws.Range("E45").Value = '  +39.81%  '

# Row 46
$ws.Range("B46").Value = 'Filecoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.01'
$ws.Range("E46").Value = '  +4.45%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0442'
$ws.Range("E47").Value = '  +5.98%  '

# Row 48
$ws.Range("E48").Value = '  +0.70%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '23.62'
$ws.Range("E49").Value = '  +0.87%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.58'
$ws.Range("E50").Value = '  +4.63%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.28'
$ws.Range("E51").Value = '  +7.32%  '

Write-Host "cryptos.xlsx updated: 93 cell(s) changed across rows 2-51"
